# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Price cells in column D are stored as plain text (values like "1.010" /
# "220.59" are NOT numbers - they're formatted strings straight from the
# scraper), so a handful of them would otherwise be auto-coerced into
# numeric values by Excel (dropping trailing zeros, e.g. "220.50" -> 220.5).
# A leading apostrophe forces those ambiguous ones to stay text, matching
# how the source file stores them; unambiguous ones (containing more than
# one '.' or special characters) are left as plain literals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.357.92"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.671.85"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'220.50"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").Value = "'0.5319"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'1.010"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.06372"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "'21.01"
$ws.Range("E10").Value = "  +2.25%  "
$ws.Range("D11").Value = "'0.07841"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "1.673.79"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "1.900.77"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "'0.5623"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").Value = "0.0₅8134"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "'66.05"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "26.370.33"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "'201.24"
$ws.Range("E21").Value = "  +4.61%  "
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").Value = "'6.073"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "'1.011"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "'0.1219"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'7.264"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("D29").Value = "'1.519"
$ws.Range("E29").Value = "  +3.15%  "
$ws.Range("D30").Value = "'0.05911"
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("D32").Value = "'3.536"
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").Value = "'1.609"
$ws.Range("D35").Value = "'0.9677"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("D36").Value = "'2.832"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").Value = "'2.431"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'0.5820"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").Value = "'0.01622"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").Value = "'5.981"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").Value = "1.079.94"
$ws.Range("E41").Value = "  +3.04%  "
$ws.Range("D42").Value = "'0.8612"
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("D44").Value = "'103.25"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").Value = "1.810.46"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").Value = "'58.64"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").Value = "'8.092"
$ws.Range("E50").Value = "  +2.08%  "
$ws.Range("D51").Value = "'0.05151"
$ws.Range("E51").Value = "  -0.29%  "
